# Scheduled-runner refresh of market-price-derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the leve
# profit sheets. Values below are the refreshed market snapshot; a few rows
# lose their LeveProfit cells entirely where the item no longer prices out.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 157.11111
$ws.Range("I12").Value = 64
$ws.Range("J12").Value = 273.5
$ws.Range("K12").Value = 64
$ws.Range("L12").Value = 273.5
$ws.Range("M12").Value = 106
$ws.Range("N12").Value = -613.5

$ws.Range("H74").Value = 2905.85
$ws.Range("I74").Value = 2683.353
$ws.Range("J74").Value = 4166.6665
$ws.Range("K74").Value = 2683.353
$ws.Range("L74").Value = 4166.6665
$ws.Range("M74").Value = -1747.353
$ws.Range("N74").Value = -6038.6665

$ws.Range("H77").Value = 2905.85
$ws.Range("I77").Value = 2683.353
$ws.Range("J77").Value = 4166.6665
$ws.Range("K77").Value = 13416.765
$ws.Range("L77").Value = 20833.3325
$ws.Range("M77").Value = -8736.764999999999
$ws.Range("N77").Value = -30193.3325

$ws.Range("H127").Value = 788.2778
$ws.Range("J127").Value = 1431.6666
$ws.Range("L127").Value = 4294.9998
$ws.Range("N127").Value = -14214.9998

$ws.Range("H129").Value = 920.1857
$ws.Range("I129").Value = 900.5
$ws.Range("J129").Value = 923.4666999999999
$ws.Range("K129").Value = 2701.5
$ws.Range("L129").Value = 2770.4001
$ws.Range("M129").Value = 2298.5
$ws.Range("N129").Value = -12770.4001

$ws.Range("H132").Value = 7817589
$ws.Range("I132").Value = 9438303
$ws.Range("J132").Value = 8691.091
$ws.Range("K132").Value = 28314909
$ws.Range("L132").Value = 26073.273
$ws.Range("M132").Value = -28312379
$ws.Range("N132").Value = -31133.273

$ws.Range("H141").Value = 2507.7058
$ws.Range("I141").Value = 1195.5
$ws.Range("J141").Value = 14580
$ws.Range("K141").Value = 3586.5
$ws.Range("L141").Value = 43740
$ws.Range("M141").Value = 1593.5
$ws.Range("N141").Value = -54100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 11038
$ws.Range("J37").Value = 11038
$ws.Range("L37").Value = 11038
$ws.Range("N37").Value = -11584

$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H61").Value = 1181.3654
$ws.Range("I61").Value = 1058.3611
$ws.Range("J61").Value = 1458.125
$ws.Range("K61").Value = 1058.3611
$ws.Range("L61").Value = 1458.125
$ws.Range("M61").Value = -846.3611000000001
$ws.Range("N61").Value = -1882.125

$ws.Range("H74").Value = 2531.1538
$ws.Range("I74").Value = 2642.0833
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 2642.0833
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -1768.0833
$ws.Range("N74").Value = -2948

$ws.Range("H77").Value = 2531.1538
$ws.Range("I77").Value = 2642.0833
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 13210.4165
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -8842.416499999999
$ws.Range("N77").Value = -14736

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 17947278
$ws.Range("I132").Value = 24391046
$ws.Range("J132").Value = 3269807
$ws.Range("K132").Value = 73173138
$ws.Range("L132").Value = 9809421
$ws.Range("M132").Value = -73170608
$ws.Range("N132").Value = -9814481

$ws.Range("H136").Value = 1181.3654
$ws.Range("I136").Value = 1058.3611
$ws.Range("J136").Value = 1458.125
$ws.Range("K136").Value = 3175.0833
$ws.Range("L136").Value = 4374.375
$ws.Range("M136").Value = -625.0833000000002
$ws.Range("N136").Value = -9474.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 30668.572
$ws.Range("J109").Value = 30668.572
$ws.Range("L109").Value = 30668.572
$ws.Range("N109").Value = -33442.572

$ws.Range("H134").Value = 3272086.5
$ws.Range("I134").Value = 1153.6
$ws.Range("K134").Value = 3460.8
$ws.Range("M134").Value = -925.7999999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1510
$ws.Range("I31").Value = 1206.4348
$ws.Range("J31").Value = 1813.5652
$ws.Range("K31").Value = 1206.4348
$ws.Range("L31").Value = 1813.5652
$ws.Range("M31").Value = -911.4348
$ws.Range("N31").Value = -2403.5652

$ws.Range("H34").Value = 1510
$ws.Range("I34").Value = 1206.4348
$ws.Range("J34").Value = 1813.5652
$ws.Range("K34").Value = 1206.4348
$ws.Range("L34").Value = 1813.5652
$ws.Range("M34").Value = -1004.4348
$ws.Range("N34").Value = -2217.5652

$ws.Range("H58").Value = 14286938
$ws.Range("I58").Value = 19609104
$ws.Range("J58").Value = 1127.2632
$ws.Range("K58").Value = 19609104
$ws.Range("L58").Value = 1127.2632
$ws.Range("M58").Value = -19608901
$ws.Range("N58").Value = -1533.2632

$ws.Range("H103").Value = 8108.8887
$ws.Range("I103").Value = 1640
$ws.Range("J103").Value = 30750
$ws.Range("K103").Value = 1640
$ws.Range("L103").Value = 30750
$ws.Range("M103").Value = -468
$ws.Range("N103").Value = -33094

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H122").Value = 41667536
$ws.Range("I122").Value = 83333670
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 250001010
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -249998560
$ws.Range("N122").Value = -9100

$ws.Range("H132").Value = 9010206
$ws.Range("I132").Value = 1083.8235
$ws.Range("J132").Value = 16667960
$ws.Range("K132").Value = 3251.4705
$ws.Range("L132").Value = 50003880
$ws.Range("M132").Value = -721.4704999999999
$ws.Range("N132").Value = -50008940

$ws.Range("H134").Value = 1024.6274
$ws.Range("I134").Value = 794.1539
$ws.Range("J134").Value = 1773.6666
$ws.Range("K134").Value = 2382.4617
$ws.Range("L134").Value = 5320.9998
$ws.Range("M134").Value = 152.5383000000002
$ws.Range("N134").Value = -10390.9998

$ws.Range("H136").Value = 14286938
$ws.Range("I136").Value = 19609104
$ws.Range("J136").Value = 1127.2632
$ws.Range("K136").Value = 58827312
$ws.Range("L136").Value = 3381.7896
$ws.Range("M136").Value = -58824762
$ws.Range("N136").Value = -8481.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11787914
$ws.Range("I131").Value = 83335000
$ws.Range("J131").Value = 5498060.5
$ws.Range("K131").Value = 250005000
$ws.Range("L131").Value = 16494181.5
$ws.Range("M131").Value = -249999960
$ws.Range("N131").Value = -16504261.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6193.316
$ws.Range("I132").Value = 4874.4546
$ws.Range("J132").Value = 8006.75
$ws.Range("K132").Value = 14623.3638
$ws.Range("L132").Value = 24020.25
$ws.Range("M132").Value = -12093.3638
$ws.Range("N132").Value = -29080.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1613.8889
$ws.Range("J22").Value = 2070.9167
$ws.Range("L22").Value = 2070.9167
$ws.Range("N22").Value = -2660.9167

$ws.Range("H27").Value = 1613.8889
$ws.Range("J27").Value = 2070.9167
$ws.Range("L27").Value = 2070.9167
$ws.Range("N27").Value = -2284.9167

$ws.Range("H46").Value = 1248.375
$ws.Range("I46").Value = 644.55554
$ws.Range("K46").Value = 644.55554
$ws.Range("M46").Value = -456.55554

$ws.Range("H55").Value = 5778.4443
$ws.Range("I55").Value = 11334.556
$ws.Range("J55").Value = 222.33333
$ws.Range("K55").Value = 11334.556
$ws.Range("L55").Value = 222.33333
$ws.Range("M55").Value = -11161.556
$ws.Range("N55").Value = -568.3333299999999

$ws.Range("H132").Value = 13517800
$ws.Range("I132").Value = 20409526
$ws.Range("J132").Value = 10014.52
$ws.Range("K132").Value = 61228578
$ws.Range("L132").Value = 30043.56
$ws.Range("M132").Value = -61226048
$ws.Range("N132").Value = -35103.56

$ws.Range("H136").Value = 24317872
$ws.Range("I136").Value = 3403206.2
$ws.Range("J136").Value = 200001070
$ws.Range("K136").Value = 10209618.6
$ws.Range("L136").Value = 600003210
$ws.Range("M136").Value = -10207068.6
$ws.Range("N136").Value = -600008310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 899.0714
$ws.Range("I81").Value = 814.2308
$ws.Range("J81").Value = 2002
$ws.Range("K81").Value = 1628.4616
$ws.Range("L81").Value = 4004
$ws.Range("M81").Value = -567.4616000000001
$ws.Range("N81").Value = -6126

$ws.Range("H84").Value = 899.0714
$ws.Range("I84").Value = 814.2308
$ws.Range("J84").Value = 2002
$ws.Range("K84").Value = 8142.308000000001
$ws.Range("L84").Value = 20020
$ws.Range("M84").Value = -2838.308000000001
$ws.Range("N84").Value = -30628

$ws.Range("H122").Value = 44658.332
$ws.Range("I122").Value = 168666.67
$ws.Range("J122").Value = 3322.2222
$ws.Range("K122").Value = 506000.01
$ws.Range("L122").Value = 9966.6666
$ws.Range("M122").Value = -503550.01
$ws.Range("N122").Value = -14866.6666

$ws.Range("H132").Value = 1214.69
$ws.Range("I132").Value = 779.6081
$ws.Range("J132").Value = 2453
$ws.Range("K132").Value = 2338.8243
$ws.Range("L132").Value = 7359
$ws.Range("M132").Value = 191.1756999999998
$ws.Range("N132").Value = -12419
